# Auto-generated edit script applying scheduled-runner price/profit updates
# across the Seraph_Profits workbook (per-sheet leve profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 192.875
$ws.Range("I4").Value = 183.6
$ws.Range("K4").Value = 183.6
$ws.Range("M4").Value = -69.59999999999999
$ws.Range("H33").Value = 185.58333
$ws.Range("I33").Value = 170.63637
$ws.Range("K33").Value = 170.63637
$ws.Range("M33").Value = 58.36363
$ws.Range("H69").Value = 11249.5
$ws.Range("J69").Value = 11249.5
$ws.Range("L69").Value = 33748.5
$ws.Range("N69").Value = -35496.5
$ws.Range("H72").Value = 11249.5
$ws.Range("J72").Value = 11249.5
$ws.Range("L72").Value = 101245.5
$ws.Range("N72").Value = -109981.5
$ws.Range("H111").Value = 1933.3334
$ws.Range("I111").Value = 1933.3334
$ws.Range("K111").Value = 5800.0002
$ws.Range("M111").Value = -2733.0002
$ws.Range("H113").Value = 3667.6667
$ws.Range("I113").Value = 3752.5
$ws.Range("J113").Value = 3498
$ws.Range("K113").Value = 3752.5
$ws.Range("L113").Value = 3498
$ws.Range("M113").Value = -498.5
$ws.Range("N113").Value = -10006
$ws.Range("H118").Value = 449.83334
$ws.Range("I118").Value = 449.83334
$ws.Range("K118").Value = 1349.50002
$ws.Range("M118").Value = 307.4999800000001
$ws.Range("H129").Value = 3103.4
$ws.Range("J129").Value = 3657.7778
$ws.Range("L129").Value = 10973.3334
$ws.Range("N129").Value = -20973.3334
$ws.Range("H135").Value = 1262.5625
$ws.Range("I135").Value = 918.2727
$ws.Range("K135").Value = 8264.454299999999
$ws.Range("M135").Value = -5729.454299999999
$ws.Range("H137").Value = 1546.7693
$ws.Range("I137").Value = 1355.875
$ws.Range("J137").Value = 1852.2
$ws.Range("K137").Value = 4067.625
$ws.Range("L137").Value = 5556.6
$ws.Range("M137").Value = -1517.625
$ws.Range("N137").Value = -10656.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 952.6923
$ws.Range("I74").Value = 764.63635
$ws.Range("K74").Value = 764.63635
$ws.Range("M74").Value = 109.36365
$ws.Range("H77").Value = 952.6923
$ws.Range("I77").Value = 764.63635
$ws.Range("K77").Value = 3823.18175
$ws.Range("M77").Value = 544.8182500000003
$ws.Range("H97").Value = 460.3158
$ws.Range("I97").Value = 482.2353
$ws.Range("J97").Value = 274
$ws.Range("K97").Value = 482.2353
$ws.Range("L97").Value = 274
$ws.Range("M97").Value = 13.7647
$ws.Range("N97").Value = -1266
$ws.Range("H122").Value = 511366.1
$ws.Range("I122").Value = 631712.7
$ws.Range("K122").Value = 1895138.1
$ws.Range("M122").Value = -1892688.1
$ws.Range("H132").Value = 1895.125
$ws.Range("I132").Value = 1878.7142
$ws.Range("J132").Value = 2010
$ws.Range("K132").Value = 5636.142599999999
$ws.Range("L132").Value = 6030
$ws.Range("M132").Value = -3106.142599999999
$ws.Range("N132").Value = -11090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2152.7222
$ws.Range("I94").Value = 1131.8334
$ws.Range("K94").Value = 1131.8334
$ws.Range("M94").Value = -680.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1532.7826
$ws.Range("I19").Value = 583.5714
$ws.Range("K19").Value = 583.5714
$ws.Range("M19").Value = -413.5714
$ws.Range("H24").Value = 1532.7826
$ws.Range("I24").Value = 583.5714
$ws.Range("K24").Value = 583.5714
$ws.Range("M24").Value = -413.5714
$ws.Range("H134").Value = 2564.1667
$ws.Range("I134").Value = 1982
$ws.Range("K134").Value = 5946
$ws.Range("M134").Value = -3411

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 6970441.5
$ws.Range("J32").Value = 6970441.5
$ws.Range("L32").Value = 20911324.5
$ws.Range("N32").Value = -20911890.5
$ws.Range("H38").Value = 124.166664
$ws.Range("I38").Value = 76.5
$ws.Range("J38").Value = 219.5
$ws.Range("K38").Value = 229.5
$ws.Range("L38").Value = 658.5
$ws.Range("M38").Value = 117.5
$ws.Range("N38").Value = -1352.5
$ws.Range("H70").Value = 537.3333
$ws.Range("I70").Value = 537.3333
$ws.Range("K70").Value = 1611.9999
$ws.Range("M70").Value = -1296.9999
$ws.Range("H73").Value = 537.3333
$ws.Range("I73").Value = 537.3333
$ws.Range("K73").Value = 1611.9999
$ws.Range("M73").Value = -519.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = ""
$ws.Range("H102").Value = 2180.3572
$ws.Range("I102").Value = 1493.125
$ws.Range("K102").Value = 1493.125
$ws.Range("M102").Value = 128.875
$ws.Range("H132").Value = 1456.4286
$ws.Range("I132").Value = 1199.1666
$ws.Range("K132").Value = 3597.4998
$ws.Range("M132").Value = -1067.4998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8388
$ws.Range("I7").Value = 9124
$ws.Range("K7").Value = 9124
$ws.Range("M7").Value = -9012
$ws.Range("H16").Value = 704.17645
$ws.Range("I16").Value = 736.8
$ws.Range("J16").Value = 459.5
$ws.Range("K16").Value = 736.8
$ws.Range("L16").Value = 459.5
$ws.Range("M16").Value = -566.8
$ws.Range("N16").Value = -799.5
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10686
$ws.Range("H40").Value = 1951
$ws.Range("I40").Value = 1951
$ws.Range("K40").Value = 1951
$ws.Range("M40").Value = -1815
$ws.Range("H126").Value = 8388
$ws.Range("I126").Value = 9124
$ws.Range("K126").Value = 27372
$ws.Range("M126").Value = -24902

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4263.3184
$ws.Range("J62").Value = 4524.9375
$ws.Range("L62").Value = 4524.9375
$ws.Range("N62").Value = -5772.9375
$ws.Range("H65").Value = 4263.3184
$ws.Range("J65").Value = 4524.9375
$ws.Range("L65").Value = 22624.6875
$ws.Range("N65").Value = -28864.6875
$ws.Range("H122").Value = 1895.6364
$ws.Range("I122").Value = 1926.75
$ws.Range("J122").Value = 1812.6666
$ws.Range("K122").Value = 5780.25
$ws.Range("L122").Value = 5437.9998
$ws.Range("M122").Value = -3330.25
$ws.Range("N122").Value = -10337.9998
$ws.Range("H132").Value = 5517.1665
$ws.Range("I132").Value = 5517.1665
$ws.Range("K132").Value = 16551.4995
$ws.Range("M132").Value = -14021.4995
